# Add a new leading "index_symbol" column (value "XYZ") to every data sheet
# that doesn't already carry a symbol column, mirroring the "yh_symbol" /
# "ms_symbol" sheets so the field can be used as part of a unique index when
# the data is stored in SQL / MongoDB.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "yh_currency",
    "ms_currency",
    "yh_esgScores",
    "earningsChart_quarterly",
    "financialsChart_yearly",
    "yh_indexTrend_estimates",
    "yh_assetProfile",
    "yh_assetProfile_companyOfficers",
    "yh_ohlcv_1d"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Shift all existing columns one to the right, freeing up column A.
    $ws.Columns.Item(1).Insert()

    # Copy the (now shifted) former-first-column header's formatting onto
    # the new A1 header cell so it keeps the same bold/centered/bordered
    # header style, then fill in the header text + values.
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("A1").Value = "index_symbol"

    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = "XYZ"
    }
}
